$d = $word.ActiveDocument

# 1) Rebuild paragraph 1's runs to match the proofed/spell-checked version,
#    using InsertXML on a range that stops just before the paragraph mark so
#    the paragraph's own identity/formatting is preserved.
$p1 = $d.Paragraphs(1)
$fullRng = $p1.Range
$rng = $d.Range($fullRng.Start, $fullRng.End - 1)

$body = @'
<w:p><w:r><w:t xml:space="preserve">2017 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>itibari</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ile</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> android </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>için</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>en</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>çok</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>kullanılan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dil</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Kotlin </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>oldu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. 2017’den </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>önce</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Java </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>popülerdi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> 2017’den </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>beri</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Android’in</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ana </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>odağı</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Kotlin’de</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>yani</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>
'@

$pkg = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $body + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($pkg) | Out-Null

# 2) Insert a new empty paragraph right before the existing trailing empty
#    paragraph (so the doc goes from 2 paragraphs to 3).
$p2 = $d.Paragraphs(2)
$insertRng = $p2.Range
$insertRng.Collapse(1)

$emptyPkg = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertRng.InsertXML($emptyPkg) | Out-Null
